$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    $cell.Value = "<br> " + $old + " <br />"
}
